# daily auto push: 2026-01-08 13:50 UTC
# A new log entry for 2026/01/08 20:00 needs to be inserted in date order,
# ahead of the existing "2026/12/29" run. Insert a fresh row at 582 (pushing
# every following row down by one) and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 582 (and everything after it) down by one row.
$ws.Rows.Item(582).Insert()

# Column A holds plain date-text like "2025/01/01", not a real Excel date.
# Force text formatting before assigning so it isn't auto-parsed into a
# date serial, then clear the formatting back to the sheet's default style
# (the data rows otherwise carry no explicit style) while keeping the
# stored value as literal text.
$ws.Range("A582").NumberFormat = "@"
$ws.Range("A582").Value = "2026/01/08"
$ws.Range("A582").ClearFormats()

$ws.Range("B582").Value = "木"
$ws.Range("C582").Value = 20
$ws.Range("D582").Value = 25
